# "forgot some reaper files"
# Update the Status column (F) on the Asset List sheet:
#   - Rows that were "Recorded" (R) and are actually finished recording/implementing
#     in Reaper now move to "Implemented" (I), highlighted with a light green fill.
#   - Rows that were "Not Started" (NS) but have since been recorded move to "R".
#   - A couple of "Not Started" rows skipped straight to "Implemented" (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose Status becomes "I" (Implemented) - gets a light green highlight fill
$implementedRows = @(2,3,9,10,12,13,14,16,17,18,31,33)

# Rows whose Status becomes "R" (Recorded) - no fill change
$recordedRows = @(20,21,22,23,25,26,27,29,30,32)

# Theme "Green, Accent 6, Lighter 60%" (#C5E0B4) expressed as an OLE BGR value
$highlightColor = 11854021

foreach ($r in $implementedRows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "I"
    $cell.Interior.Color = $highlightColor
}

foreach ($r in $recordedRows) {
    $cell = $ws.Cells.Item($r, 6)
    $cell.Value = "R"
}

# Restore the view roughly where the author left off scrolling/selecting
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("F28").Select() | Out-Null
